# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" (column D, row 3) on both the
# zh-cn and de-de status sheets to reflect the new handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-03-10 02:20:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-03-10 02:20:23"
